$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C (Forest area (ha)) from 21.71 to 23.71
# (Excel snaps ColumnWidth to the nearest renderable pixel width for the
# sheet's default font, so we feed the value whose snapped result lands
# closest to the target stored width of 23.7109375 characters.)
$ws.Columns.Item(3).ColumnWidth = 22.8333333333333

$data = @{
    2  = @{ A = 255;  B = 9210;   C = 30.70000076293945; H = 12;  I = -22039; J = 214 }
    3  = @{ A = 62;   B = 646;    C = 14.04347801208496; H = -26; I = -74862; J = 100 }
    4  = @{ A = 53;   B = 898;    C = 7.126984119415283; H = 99;  I = 92098;  J = -70 }
    5  = @{ A = 166;  B = 9564;   C = 13.28333377838135; H = 124; I = 100146; J = 323 }
    6  = @{ A = 288;  B = 34271;  C = 35.73618316650391; H = 94;  I = 82383;  J = 301 }
    7  = @{ A = 217;  B = 58652;  C = 64.3114013671875;  H = 192; I = 205232; J = 379 }
    8  = @{ A = 131;  B = 10614;  C = 23.43046379089355; H = 52;  I = 8872;   J = 207 }
    9  = @{ A = 142;  B = 26795;  C = 79.74702453613281; H = 46;  I = 25643;  J = 17  }
    10 = @{ A = 99;   B = 5100;   C = 21.51898765563965; H = 22;  I = -6824;  J = 86  }
    11 = @{ A = 131;  B = 8732;   C = 16.50661659240723; H = 59;  I = 14230;  J = 219 }
    12 = @{ A = 260;  B = 61888;  C = 149.1277160644531; H = 117; I = 75887;  J = 33  }
    13 = @{ A = 228;  B = 93465;  C = 113.2909088134766; H = 17;  I = -15419; J = 134 }
    14 = @{ A = 1067; B = 786198; C = 636.5975952148438; H = 567; I = 668907; J = 548 }
    15 = @{ A = 387;  B = 134647; C = 260.4390563964844; H = 170; I = 97885;  J = 339 }
    16 = @{ A = 64;   B = 1775;   C = 18.88297843933105; H = 87;  I = 20496;  J = 236 }
    17 = @{ A = 85;   B = 2338;   C = 19.1639347076416;  H = -6;  I = -42198; J = 175 }
    18 = @{ A = 118;  B = 18388;  C = 75.36065673828125; H = 66;  I = 4821;   J = 273 }
    19 = @{ A = 259;  B = 81724;  C = 213.3785858154297; H = 17;  I = -12387; J = 189 }
    20 = @{ A = 44;   B = 12674;  C = 22.7132625579834;  H = 78;  I = 31670;  J = 81  }
    21 = @{ A = 57;   B = 9198;   C = 16.9392261505127;  H = 74;  I = 50278;  J = 272 }
    22 = @{ A = 339;  B = 185155; C = 261.9937133789062; H = 234; I = 247112; J = 395 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
}
